# Module 2 Week 1 Example
# Adds a duplicate-transaction check in column F of Sheet1:
#   F2:F11 = COUNTIFS(B:B,B,C:C,C,D:D,D) > 1   (TRUE if the Name/Date/Amount combo repeats)
# and leaves the selection on F4, matching the authored workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# First row of the duplicate-check column (single formula, not part of the shared group).
$ws.Range("F2").Formula = "=COUNTIFS(B2:B11,B2,C2:C11,C2,D2:D11,D2)>1"

# Remaining rows share one relative formula, mirroring the existing pattern used in column G.
$ws.Range("F3:F11").Formula = "=COUNTIFS(B3:B12,B3,C3:C12,C3,D3:D12,D3)>1"

# Leave the active cell on F4 as in the saved workbook.
$ws.Range("F4").Select() | Out-Null
